$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto Price (D) and 1h Volume-change (E) columns with the latest
# scrape. Values such as "1.00" / "0.205" would otherwise be auto-coerced to
# plain numbers (losing the trailing zero / exact text) by Excel's normal
# "looks like a number" detection, so those cells are briefly marked as Text
# (NumberFormat "@") while the value is written, then restored to the default
# "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = '72.022.55'
$ws.Range("E2").Value = '  +3.93%  '

$ws.Range("D3").Value = '3.701.14'
$ws.Range("E3").Value = '  +9.65%  '

$ws.Range("E4").Value = '  -0.10%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '589.38'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.66%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '181.12'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.68%  '

$ws.Range("D7").Value = '3.698.38'
$ws.Range("E7").Value = '  +9.77%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.625'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +5.38%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.01%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.205'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +3.75%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.614'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +4.74%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '50.08'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +3.58%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.0000289'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.07%  '

$ws.Range("D14").Value = '4.285.04'
$ws.Range("E14").Value = '  +9.36%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '687.67'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.08%  '

$ws.Range("E16").Value = '  +4.82%  '

$ws.Range("D17").Value = '3.731.75'
$ws.Range("E17").Value = '  +10.39%  '

$ws.Range("D18").Value = '72.102.15'
$ws.Range("E18").Value = '  +3.87%  '

$ws.Range("E19").Value = '  +2.15%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '18.33'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +3.30%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '11.70'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +4.17%  '

$ws.Range("E22").Value = '  +4.42%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.19'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +15.92%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '17.86'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +3.88%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '103.90'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +2.52%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '4.05'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +4.50%  '

$ws.Range("E27").Value = '  +6.45%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.32'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +6.51%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '35.58'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +6.35%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '9.25'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +6.28%  '

$ws.Range("E31").Value = '  +7.78%  '

$ws.Range("E32").Value = '  +10.89%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '579.34'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +5.11%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '11.35'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +3.02%  '

$ws.Range("E35").Value = '  +4.19%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '60.27'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.22%  '

$ws.Range("D37").Value = '3.752.16'
$ws.Range("E37").Value = '  +4.26%  '

$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("E39").Value = '  +3.05%  '

$ws.Range("D40").Value = '0.0₃0779'
$ws.Range("E40").Value = '  +7.41%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '35.75'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.47%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '3.48'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +5.29%  '

$ws.Range("E43").Value = '  +3.33%  '

$ws.Range("E44").Value = '  +8.92%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.349'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +4.12%  '

$ws.Range("E46").Value = '  +0.03%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '2.85'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +7.45%  '

$ws.Range("E48").Value = '  +4.10%  '

$ws.Range("E49").Value = '  +4.85%  '

$ws.Range("E50").Value = '  -0.30%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '134.09'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +3.97%  '
